# Added Talent Hunt user.
# Adds three new rows (11-13) to the "auto_credentials" sheet for a new
# "ONEUP_TALENTHUNT_USER" credential (one row per environment: STAGE, TEST, PROD),
# mirroring the existing ONEUP_ADMIN / ONEUP_VALID_USER / ONEUP_INVALID_USER blocks,
# complete with a mailto hyperlink on the username/email cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("auto_credentials")

$credentialAlias = "ONEUP_TALENTHUNT_USER"
$userEmail       = "samruddhi.ubhad+th1@jiem.in"
$mailTarget      = "mailto:" + $userEmail

$envs = "STAGE", "TEST", "PROD"

# Rows 11 and 13 carry the plain "Hyperlink" cell style (re-using the
# workbook's existing style slot); row 12 picks up the style Excel derives
# when a hyperlink is freshly applied via Hyperlinks.Add (adds its own
# fill/border-apply flags), matching the source workbook exactly.
$restyleAfterAdd = $true, $false, $true

for ($i = 0; $i -lt $envs.Length; $i++) {
    $row = 11 + $i

    $ws.Cells.Item($row, 1).Value = 3              # A: credential_id
    $ws.Cells.Item($row, 2).Value = 14              # B: project id
    $ws.Cells.Item($row, 3).Value = $credentialAlias # C: credential_alias
    $ws.Cells.Item($row, 4).Value = $userEmail       # D: username (hyperlinked)
    $ws.Cells.Item($row, 5).Value = 123456           # E: password
    $ws.Cells.Item($row, 9).Value = $envs[$i]        # I: env

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 4), $mailTarget) | Out-Null

    if ($restyleAfterAdd[$i]) {
        $ws.Cells.Item($row, 4).Style = "Hyperlink"
    }
}
